$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 303.84128
$ws.Range("J17").Value = 303.84128
$ws.Range("L17").Value = 911.52384
$ws.Range("N17").Value = -1247.52384

$ws.Range("H41").Value = 3339.0588
$ws.Range("J41").Value = 6303.294
$ws.Range("L41").Value = 6303.294
$ws.Range("N41").Value = -7183.294

$ws.Range("H62").Value = 9826.655000000001
$ws.Range("I62").Value = 13504.412
$ws.Range("J62").Value = 4616.5
$ws.Range("K62").Value = 13504.412
$ws.Range("L62").Value = 4616.5
$ws.Range("M62").Value = -12880.412
$ws.Range("N62").Value = -5864.5

$ws.Range("H64").Value = 4986.4546
$ws.Range("I64").Value = 4715.3076
$ws.Range("J64").Value = 5378.1113
$ws.Range("K64").Value = 4715.3076
$ws.Range("L64").Value = 5378.1113
$ws.Range("M64").Value = -4467.3076
$ws.Range("N64").Value = -5874.1113

$ws.Range("H65").Value = 9826.655000000001
$ws.Range("I65").Value = 13504.412
$ws.Range("J65").Value = 4616.5
$ws.Range("K65").Value = 67522.06
$ws.Range("L65").Value = 23082.5
$ws.Range("M65").Value = -64402.06
$ws.Range("N65").Value = -29322.5

$ws.Range("H67").Value = 4986.4546
$ws.Range("I67").Value = 4715.3076
$ws.Range("J67").Value = 5378.1113
$ws.Range("K67").Value = 4715.3076
$ws.Range("L67").Value = 5378.1113
$ws.Range("M67").Value = -3857.3076
$ws.Range("N67").Value = -7094.1113

$ws.Range("H116").Value = 41003.75
$ws.Range("I116").Value = 65200.41
$ws.Range("J116").Value = 3608.9092
$ws.Range("K116").Value = 65200.41
$ws.Range("L116").Value = 3608.9092
$ws.Range("M116").Value = -61758.41
$ws.Range("N116").Value = -10492.9092

$ws.Range("H132").Value = 2767.6345
$ws.Range("I132").Value = 1599.5555
$ws.Range("J132").Value = 10276.714
$ws.Range("K132").Value = 4798.666499999999
$ws.Range("L132").Value = 30830.142
$ws.Range("M132").Value = -2268.666499999999
$ws.Range("N132").Value = -35890.142

$ws.Range("H137").Value = 1999.2363
$ws.Range("I137").Value = 1995.0256
$ws.Range("J137").Value = 2009.5
$ws.Range("K137").Value = 5985.0768
$ws.Range("L137").Value = 6028.5
$ws.Range("M137").Value = -3435.0768
$ws.Range("N137").Value = -11128.5


# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 188058.52
$ws.Range("I74").Value = 251061.5
$ws.Range("J74").Value = 55420.633
$ws.Range("K74").Value = 251061.5
$ws.Range("L74").Value = 55420.633
$ws.Range("M74").Value = -250187.5
$ws.Range("N74").Value = -57168.633

$ws.Range("H77").Value = 188058.52
$ws.Range("I77").Value = 251061.5
$ws.Range("J77").Value = 55420.633
$ws.Range("K77").Value = 1255307.5
$ws.Range("L77").Value = 277103.165
$ws.Range("M77").Value = -1250939.5
$ws.Range("N77").Value = -285839.165

$ws.Range("H97").Value = 1116.08
$ws.Range("I97").Value = 1308.125
$ws.Range("J97").Value = 774.6667
$ws.Range("K97").Value = 1308.125
$ws.Range("L97").Value = 774.6667
$ws.Range("M97").Value = -812.125
$ws.Range("N97").Value = -1766.6667

$ws.Range("H110").Value = 1396.5454
$ws.Range("I110").Value = 1573.5555
$ws.Range("J110").Value = 600
$ws.Range("K110").Value = 1573.5555
$ws.Range("L110").Value = 600
$ws.Range("M110").Value = 471.4445000000001
$ws.Range("N110").Value = -4690

$ws.Range("H122").Value = 3916.6924
$ws.Range("I122").Value = 3597.7646
$ws.Range("K122").Value = 10793.2938
$ws.Range("M122").Value = -8343.293799999999

$ws.Range("H132").Value = 28174.072
$ws.Range("I132").Value = 47369.348
$ws.Range("J132").Value = 3646.7778
$ws.Range("K132").Value = 142108.044
$ws.Range("L132").Value = 10940.3334
$ws.Range("M132").Value = -139578.044
$ws.Range("N132").Value = -16000.3334


# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 30000
$ws.Range("J9").Value = 30000
$ws.Range("L9").Value = 30000
$ws.Range("N9").Value = -30336


# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2343.9856
$ws.Range("I31").Value = 1802.25
$ws.Range("K31").Value = 1802.25
$ws.Range("M31").Value = -1507.25

$ws.Range("H34").Value = 2343.9856
$ws.Range("I34").Value = 1802.25
$ws.Range("K34").Value = 1802.25
$ws.Range("M34").Value = -1600.25

$ws.Range("H105").Value = 976.0714
$ws.Range("I105").Value = 720.55554
$ws.Range("J105").Value = 1436
$ws.Range("K105").Value = 720.55554
$ws.Range("L105").Value = 1436
$ws.Range("M105").Value = 1026.44446
$ws.Range("N105").Value = -4930

$ws.Range("H134").Value = 1829.8966
$ws.Range("I134").Value = 1187.0555
$ws.Range("J134").Value = 2881.818
$ws.Range("K134").Value = 3561.1665
$ws.Range("L134").Value = 8645.454000000002
$ws.Range("M134").Value = -1026.1665
$ws.Range("N134").Value = -13715.454


# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3066.5625
$ws.Range("I129").Value = 5538.3335
$ws.Range("J129").Value = 1583.5
$ws.Range("K129").Value = 16615.0005
$ws.Range("L129").Value = 4750.5
$ws.Range("M129").Value = -11615.0005
$ws.Range("N129").Value = -14750.5

$ws.Range("H132").Value = 3724.1516
$ws.Range("I132").Value = 2187.7334
$ws.Range("J132").Value = 5004.5
$ws.Range("K132").Value = 19689.6006
$ws.Range("L132").Value = 45040.5
$ws.Range("M132").Value = -17159.6006
$ws.Range("N132").Value = -50100.5


# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4921.8696
$ws.Range("I102").Value = 3825.7693
$ws.Range("J102").Value = 6346.8
$ws.Range("K102").Value = 3825.7693
$ws.Range("L102").Value = 6346.8
$ws.Range("M102").Value = -2203.7693
$ws.Range("N102").Value = -9590.799999999999

$ws.Range("H122").Value = 1491.091
$ws.Range("I122").Value = 1618.5454
$ws.Range("J122").Value = 1363.6364
$ws.Range("K122").Value = 4855.6362
$ws.Range("L122").Value = 4090.9092
$ws.Range("M122").Value = -2405.6362
$ws.Range("N122").Value = -8990.9092

$ws.Range("H126").Value = 3245.3447
$ws.Range("I126").Value = 3050.9092
$ws.Range("J126").Value = 3364.1667
$ws.Range("K126").Value = 9152.7276
$ws.Range("L126").Value = 10092.5001
$ws.Range("M126").Value = -6682.7276
$ws.Range("N126").Value = -15032.5001

$ws.Range("H132").Value = 3688.814
$ws.Range("I132").Value = 3658.6155
$ws.Range("J132").Value = 3735
$ws.Range("K132").Value = 10975.8465
$ws.Range("L132").Value = 11205
$ws.Range("M132").Value = -8445.8465
$ws.Range("N132").Value = -16265


# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9838.275
$ws.Range("I132").Value = 3472.182
$ws.Range("K132").Value = 10416.546
$ws.Range("M132").Value = -7886.545999999998


# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2472.7273
$ws.Range("I81").Value = 1157.1428
$ws.Range("J81").Value = 4775
$ws.Range("K81").Value = 2314.2856
$ws.Range("L81").Value = 9550
$ws.Range("M81").Value = -1253.2856
$ws.Range("N81").Value = -11672

$ws.Range("H84").Value = 2472.7273
$ws.Range("I84").Value = 1157.1428
$ws.Range("J84").Value = 4775
$ws.Range("K84").Value = 11571.428
$ws.Range("L84").Value = 47750
$ws.Range("M84").Value = -6267.428
$ws.Range("N84").Value = -58358

$ws.Range("H132").Value = 1739.674
$ws.Range("I132").Value = 1056.0278
$ws.Range("K132").Value = 3168.0834
$ws.Range("M132").Value = -638.0834000000004

